# DB Schema and Initial Dump
# - Rename the "Users" sheet's "create_account" permission column (D) to "register"
# - Rename the "Users" sheet's "create_special_account" permission column (J) to "create_account"
# - Switch focus to the "Users" sheet, with column D selected (active cell D1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("D1").Value = "register"
$ws.Range("J1").Value = "create_account"

$ws.Activate()
$ws.Columns.Item(4).Select()
